# Insert a new data row at row 316 (pushing existing rows 316-445 down to 317-446)
# and populate the new row with the latest weekly price observation for Coliflor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 316:445 down to 317:446, creating a blank row at 316
$ws.Rows.Item(316).Insert()

# Populate the new row 316 with the new record's data
$ws.Cells.Item(316, 1).Value = 10
$ws.Cells.Item(316, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(316, 3).Value = "La Araucanía"
$ws.Cells.Item(316, 4).Value = 44795
$ws.Cells.Item(316, 5).Value = 9
$ws.Cells.Item(316, 6).Value = 100112008
$ws.Cells.Item(316, 7).Value = "Coliflor"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 5000
$ws.Cells.Item(316, 11).Value = 1300
$ws.Cells.Item(316, 12).Value = 1400
$ws.Cells.Item(316, 13).Value = 1340
$ws.Cells.Item(316, 14).Value = "$/unidad"
$ws.Cells.Item(316, 15).Value = "Región Metropolitana"
$ws.Cells.Item(316, 16).Value = 1340
$ws.Cells.Item(316, 17).Value = 1
$ws.Cells.Item(316, 18).Value = "Hortaliza"
